# Updates evaluation metric values (wa_mpjpe, w_mpjpe, rte, erve, ate, ate_s)
# for the affected sequences following the fix for scale estimation
# (commit: "add mono mast3r & fix scale estimation").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 108.6829299926758
$ws.Range("H2").Value = 317.6914978027344
$ws.Range("I2").Value = 3.086787939071655
$ws.Range("J2").Value = 10.21835327148438
$ws.Range("K2").Value = 0.430992442154261
$ws.Range("L2").Value = 1.990642279791056
$ws.Range("G3").Value = 99.13406372070312
$ws.Range("H3").Value = 216.0992126464844
$ws.Range("I3").Value = 6.81588077545166
$ws.Range("J3").Value = 6.334219932556152
$ws.Range("K3").Value = 0.08980409620695536
$ws.Range("L3").Value = 1.259197856155351
$ws.Range("G4").Value = 165.6681060791016
$ws.Range("H4").Value = 740.0010986328125
$ws.Range("I4").Value = 7.536784648895264
$ws.Range("J4").Value = 16.80570793151855
$ws.Range("K4").Value = 3.225817004402545
$ws.Range("L4").Value = 6.791158183206554
$ws.Range("G6").Value = 73.77909851074219
$ws.Range("H6").Value = 273.7779846191406
$ws.Range("I6").Value = 2.889680862426758
$ws.Range("J6").Value = 9.387054443359375
$ws.Range("K6").Value = 0.4648379908907747
$ws.Range("L6").Value = 1.229267315286151
$ws.Range("G8").Value = 116.0167922973633
$ws.Range("H8").Value = 435.6024169921875
$ws.Range("I8").Value = 5.029820442199707
$ws.Range("J8").Value = 20.85674095153809
$ws.Range("K8").Value = 0.5089652009826994
$ws.Range("L8").Value = 1.983294202847733
$ws.Range("G9").Value = 166.8375091552734
$ws.Range("H9").Value = 443.2389831542969
$ws.Range("I9").Value = 4.949410438537598
$ws.Range("J9").Value = 16.43637275695801
$ws.Range("K9").Value = 1.317884681673736
$ws.Range("L9").Value = 1.501406740404275
$ws.Range("G10").Value = 47.57673645019531
$ws.Range("H10").Value = 113.9475021362305
$ws.Range("I10").Value = 0.3832592964172363
$ws.Range("J10").Value = 8.047077178955078
$ws.Range("K10").Value = 0.08098722138961886
$ws.Range("L10").Value = 0.1102630363814392
$ws.Range("G12").Value = 57.89563751220703
$ws.Range("H12").Value = 155.9514923095703
$ws.Range("I12").Value = 0.2485483884811401
$ws.Range("J12").Value = 13.51173210144043
$ws.Range("K12").Value = 0.03533377400280525
$ws.Range("L12").Value = 0.03621065308990438
$ws.Range("G14").Value = 249.5182952880859
$ws.Range("H14").Value = 668.255615234375
$ws.Range("I14").Value = 7.366587162017822
$ws.Range("J14").Value = 14.78339099884033
$ws.Range("K14").Value = 3.614055647838236
$ws.Range("L14").Value = 4.433842816713573
$ws.Range("G15").Value = 124.0782470703125
$ws.Range("H15").Value = 227.3143310546875
$ws.Range("I15").Value = 4.94888162612915
$ws.Range("J15").Value = 7.703754425048828
$ws.Range("K15").Value = 0.1348285218486699
$ws.Range("L15").Value = 0.7558165751667402
$ws.Range("G16").Value = 79.10270690917969
$ws.Range("H16").Value = 245.5284118652344
$ws.Range("I16").Value = 0.9533370137214661
$ws.Range("J16").Value = 8.215726852416992
$ws.Range("K16").Value = 0.5461294901259565
$ws.Range("L16").Value = 0.610382169987709
$ws.Range("G17").Value = 115.2895736694336
$ws.Range("H17").Value = 317.4119262695312
$ws.Range("I17").Value = 1.091237187385559
$ws.Range("J17").Value = 13.26593208312988
$ws.Range("K17").Value = 0.3075736701257018
$ws.Range("L17").Value = 0.388819887458323
$ws.Range("G18").Value = 90.23635864257812
$ws.Range("H18").Value = 187.8910522460938
$ws.Range("I18").Value = 0.8822309970855713
$ws.Range("J18").Value = 9.339920043945312
$ws.Range("K18").Value = 0.1272548079570197
$ws.Range("L18").Value = 0.2722432278784981
$ws.Range("G19").Value = 125.7876358032227
$ws.Range("H19").Value = 521.701171875
$ws.Range("I19").Value = 5.127596855163574
$ws.Range("J19").Value = 17.47074127197266
$ws.Range("K19").Value = 0.6956230466389091
$ws.Range("L19").Value = 2.540412812436776
$ws.Range("G21").Value = 137.5004425048828
$ws.Range("H21").Value = 408.0429382324219
$ws.Range("I21").Value = 2.171093463897705
$ws.Range("J21").Value = 9.779004096984863
$ws.Range("K21").Value = 0.8812014897652418
$ws.Range("L21").Value = 1.230758937724136
$ws.Range("G23").Value = 77.65883636474609
$ws.Range("H23").Value = 184.5671997070312
$ws.Range("I23").Value = 1.013454437255859
$ws.Range("J23").Value = 9.039108276367188
$ws.Range("K23").Value = 0.1385385312536891
$ws.Range("L23").Value = 0.166911543036723
$ws.Range("G24").Value = 88.29719543457031
$ws.Range("H24").Value = 237.0995788574219
$ws.Range("I24").Value = 0.79860520362854
$ws.Range("J24").Value = 13.37563800811768
$ws.Range("K24").Value = 0.1596388178446454
$ws.Range("L24").Value = 0.1775355738014216
$ws.Range("G25").Value = 56.22836303710938
$ws.Range("H25").Value = 210.8364410400391
$ws.Range("I25").Value = 0.709577739238739
$ws.Range("J25").Value = 10.7540111541748
$ws.Range("K25").Value = 0.2740642153934587
$ws.Range("L25").Value = 0.6239872640883656
$ws.Range("G27").Value = 133.3283081054688
$ws.Range("H27").Value = 416.593017578125
$ws.Range("I27").Value = 4.061065196990967
$ws.Range("J27").Value = 12.60750770568848
$ws.Range("K27").Value = 0.8907501605608155
$ws.Range("L27").Value = 2.172430600232776
